$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-11, columns B-G (ME, MAE, MSE, RMSE, SE, N)
# Row 2 receives a brand-new computed data point (shifted in from the naive
# forecaster fix); rows 3-11 take on the values that previously belonged to
# rows 2-10 (i.e. the whole block shifted down by one row), and the former
# row 11 values are dropped.

$data = @{
    2  = @(-0.09914768942982544, 0.6760636891099804, 1.058759192261776, 1.02896024814459, 1.052237025233689, 19)
    3  = @(-0.1143386297721907, 0.6889840097079422, 0.9692665404913205, 0.9845133521142924, 1.006200714633619, 18)
    4  = @(-0.03231471279853988, 0.7695157101427289, 0.9504010021523444, 0.9748851225412891, 1.004336377456358, 17)
    5  = @(0.05671622629984657, 0.6321910596607474, 0.8518300650585396, 0.9229464042177854, 0.951413452211165, 16)
    6  = @(0.085928964333323, 0.5612671955956002, 0.5873845025919633, 0.7664101399328973, 0.7883078955770884, 15)
    7  = @(0.1382138394267656, 0.7476276879241753, 1.097975515794379, 1.047843268716452, 1.077897322974905, 14)
    8  = @(0.1215206328360093, 0.7819526689445239, 1.205421270947194, 1.097916786895616, 1.135726716235405, 13)
    9  = @(0.2028461735751207, 0.8472735835465873, 1.158512388347609, 1.076342133500128, 1.10405820904526, 12)
    10 = @(0.239907996146195, 0.6651045147971374, 0.5928596743010739, 0.769973814035954, 0.7673553351966808, 11)
    11 = @(0.2225839694609408, 0.8306889768957216, 1.245987321698922, 1.116238021973325, 1.152988242321246, 10)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

$wb.Save()
